# Auto-generated edit script applying the Raiden_Profits market-data refresh
# across all 8 item-category sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1818.4736
$ws.Range("I28").Value = 1058.3846
$ws.Range("J28").Value = 3465.3333
$ws.Range("K28").Value = 1058.3846
$ws.Range("L28").Value = 3465.3333
$ws.Range("M28").Value = -573.3846000000001
$ws.Range("N28").Value = -4435.3333
# Row 32
$ws.Range("H32").Value = 4998.5
$ws.Range("J32").Value = 4999
$ws.Range("L32").Value = 4999
$ws.Range("N32").Value = -5651
# Row 74
$ws.Range("H74").Value = 4422.5
$ws.Range("I74").Value = 4422.5
$ws.Range("K74").Value = 4422.5
$ws.Range("M74").Value = -3486.5
# Row 77
$ws.Range("H77").Value = 4422.5
$ws.Range("I77").Value = 4422.5
$ws.Range("K77").Value = 22112.5
$ws.Range("M77").Value = -17432.5
# Row 107
$ws.Range("H107").Value = 1244.5385
$ws.Range("I107").Value = 258.81818
$ws.Range("J107").Value = 6666
$ws.Range("K107").Value = 258.81818
$ws.Range("L107").Value = 6666
$ws.Range("M107").Value = 1661.18182
$ws.Range("N107").Value = -10506
# Row 112
$ws.Range("H112").Value = 1836.125
$ws.Range("J112").Value = 2093.4546
$ws.Range("L112").Value = 6280.3638
$ws.Range("N112").Value = -8496.363799999999
# Row 116
$ws.Range("H116").Value = 4386
$ws.Range("I116").Value = 4326.6665
$ws.Range("K116").Value = 4326.6665
$ws.Range("M116").Value = -884.6665000000003
# Row 138
$ws.Range("H138").Value = 2949.4854
$ws.Range("J138").Value = 2467.0378
$ws.Range("L138").Value = 7401.1134
$ws.Range("N138").Value = -17681.1134

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1992.0588
$ws.Range("I32").Value = 1624.4286
$ws.Range("K32").Value = 1624.4286
$ws.Range("M32").Value = -1337.4286
# Row 63
$ws.Range("H63").Value = 6640.5713
$ws.Range("I63").Value = 5621
$ws.Range("K63").Value = 5621
$ws.Range("M63").Value = -4935
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
# Row 66
$ws.Range("H66").Value = 6640.5713
$ws.Range("I66").Value = 5621
$ws.Range("K66").Value = 28105
$ws.Range("M66").Value = -24673
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
# Row 76
$ws.Range("H76").Value = 82670
$ws.Range("I76").Value = 82670
$ws.Range("K76").Value = 82670
$ws.Range("M76").Value = -82332
# Row 79
$ws.Range("H79").Value = 82670
$ws.Range("I79").Value = 82670
$ws.Range("K79").Value = 82670
$ws.Range("M79").Value = -81500
# Row 80
$ws.Range("H80").Value = 19980
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 19980
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 19980
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -21976
# Row 83
$ws.Range("H83").Value = 19980
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 19980
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 59940
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -69924
# Row 88
$ws.Range("H88").Value = 3949.5
$ws.Range("J88").Value = 3899
$ws.Range("L88").Value = 3899
$ws.Range("N88").Value = -4711
# Row 91
$ws.Range("H91").Value = 3949.5
$ws.Range("J91").Value = 3899
$ws.Range("L91").Value = 3899
$ws.Range("N91").Value = -6707
# Row 97
$ws.Range("H97").Value = 1107.1
$ws.Range("I97").Value = 507.8889
$ws.Range("K97").Value = 507.8889
$ws.Range("M97").Value = -11.88889999999998
# Row 102
$ws.Range("H102").Value = 1005
$ws.Range("I102").Value = 1005
$ws.Range("K102").Value = 1005
$ws.Range("M102").Value = 617
# Row 110
$ws.Range("H110").Value = 1692
$ws.Range("J110").Value = 1438.4
$ws.Range("L110").Value = 1438.4
$ws.Range("N110").Value = -5528.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 7700
$ws.Range("I94").Value = 900
$ws.Range("K94").Value = 900
$ws.Range("M94").Value = -449
# Row 99
$ws.Range("H99").Value = 4328.8
$ws.Range("I99").Value = 4274
$ws.Range("J99").Value = 4548
$ws.Range("K99").Value = 4274
$ws.Range("L99").Value = 4548
$ws.Range("M99").Value = -2776
$ws.Range("N99").Value = -7544

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 4067.5715
$ws.Range("I94").Value = 4494.8
$ws.Range("K94").Value = 4494.8
$ws.Range("M94").Value = -4043.8
# Row 95
$ws.Range("H95").Value = 51000
$ws.Range("J95").Value = 51000
$ws.Range("L95").Value = 51000
$ws.Range("N95").Value = -56492
# Row 107
$ws.Range("H107").Value = 678.25
$ws.Range("I107").Value = 678.25
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 678.25
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1241.75
$ws.Range("N107").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 337.6
$ws.Range("I14").Value = 337.6
$ws.Range("K14").Value = 1012.8
$ws.Range("M14").Value = -839.8000000000001
# Row 131
$ws.Range("H131").Value = 1657632.1
$ws.Range("J131").Value = 1774992.1
$ws.Range("L131").Value = 5324976.300000001
$ws.Range("N131").Value = -5335056.300000001
# Row 132
$ws.Range("H132").Value = 225
$ws.Range("I132").Value = 225
$ws.Range("K132").Value = 2025
$ws.Range("M132").Value = 505
# Row 134
$ws.Range("H134").Value = 697.6667
$ws.Range("I134").Value = 697.6667
$ws.Range("K134").Value = 2093.0001
$ws.Range("M134").Value = 2976.9999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 16503.23
$ws.Range("I80").Value = 8192.5
$ws.Range("J80").Value = 23626.715
$ws.Range("K80").Value = 8192.5
$ws.Range("L80").Value = 23626.715
$ws.Range("M80").Value = -7194.5
$ws.Range("N80").Value = -25622.715
# Row 83
$ws.Range("H83").Value = 16503.23
$ws.Range("I83").Value = 8192.5
$ws.Range("J83").Value = 23626.715
$ws.Range("K83").Value = 40962.5
$ws.Range("L83").Value = 118133.575
$ws.Range("M83").Value = -35970.5
$ws.Range("N83").Value = -128117.575
# Row 107
$ws.Range("H107").Value = 1040.6
$ws.Range("J107").Value = 1244.1428
$ws.Range("L107").Value = 1244.1428
$ws.Range("N107").Value = -5084.1428
# Row 113
$ws.Range("H113").Value = 3444.6667
$ws.Range("I113").Value = 2260.6
$ws.Range("K113").Value = 2260.6
$ws.Range("M113").Value = -90.59999999999991
# Row 136
$ws.Range("H136").Value = 50462.832
$ws.Range("J136").Value = 50462.832
$ws.Range("L136").Value = 151388.496
$ws.Range("N136").Value = -156488.496

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 339.1
$ws.Range("J16").Value = 329.5
$ws.Range("L16").Value = 329.5
$ws.Range("N16").Value = -669.5
# Row 93
$ws.Range("H93").Value = 3320
$ws.Range("I93").Value = 3320
$ws.Range("K93").Value = 3320
$ws.Range("M93").Value = -2072
# Row 100
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 92
$ws.Range("H92").Value = 32659.666
$ws.Range("J92").Value = 32659.666
$ws.Range("L92").Value = 32659.666
$ws.Range("N92").Value = -37651.666
# Row 94
$ws.Range("H94").Value = 24743.334
$ws.Range("J94").Value = 24743.334
$ws.Range("L94").Value = 24743.334
$ws.Range("N94").Value = -26545.334
# Row 100
$ws.Range("H100").Value = 410.2
$ws.Range("I100").Value = 387.875
$ws.Range("K100").Value = 775.75
$ws.Range("M100").Value = -234.75
# Row 122
$ws.Range("H122").Value = 4449.2188
$ws.Range("I122").Value = 4391.778
$ws.Range("K122").Value = 13175.334
$ws.Range("M122").Value = -10725.334
# Row 126
$ws.Range("H126").Value = 858
$ws.Range("I126").Value = 822.5
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 2467.5
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = 2.5
$ws.Range("N126").Value = -7940
